$wb = $excel.ActiveWorkbook

# --- Algorithm 2 sheet: selection only moved (E18 -> I2) ---
$ws2 = $wb.Worksheets.Item("Algorithm 2")
$ws2.Range("I2").Select() | Out-Null

# --- Algorithm 3 sheet: the "fixed c++ block implementation" data fill-in ---
$ws3 = $wb.Worksheets.Item("Algorithm 3")
$ws3.Activate()

# New row/col label header for the block-size table
$ws3.Range("A1").Value = "Block Size\Data Size"

# Widen column A so the new label fits
$ws3.Columns.Item(1).ColumnWidth = 19.28515625

# Newly measured results for block sizes 128 / 256 / 512
$ws3.Range("B2").Value = 26.37
$ws3.Range("C2").Value = 87.391
$ws3.Range("D2").Value = 209.766

$ws3.Range("B3").Value = 24.9
$ws3.Range("C3").Value = 87.281
$ws3.Range("D3").Value = 211.234

$ws3.Range("B4").Value = 20.656
$ws3.Range("C4").Value = 75.797
$ws3.Range("D4").Value = 190.406
$ws3.Range("E4").Value = 380.344

# Empty styled placeholder cell below the table (matches underline style used elsewhere)
$ws3.Range("D5").Font.Underline = $true

# Page setup: portrait orientation
$ws3.PageSetup.Orientation = 1

# Restore Algorithm 3's own selection state last (it is the active/tabbed sheet)
$ws3.Range("E5").Select() | Out-Null
